$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ups")

# Rename column header D1 from "IsAvailable" to "Availability"
$ws.Range("D1").Value = "Availability"

# Column D used to hold boolean TRUE values; fix the data so it stores
# plain numeric 0 for every product row (2-32) instead of the boolean.
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
}

# Restore the active selection to E19
$ws.Range("E19").Select()
